$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new columns (AD, AE, AF) for the team's season record:
# Wins, Losses, Ties. Copy the existing header style (from AC1, the
# last header cell) so the new headers look consistent with the rest
# of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill every data row (2-49) with the team's season record: 90 wins,
# 72 losses, 0 ties.
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}
